$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===================================================================
# Step 1: stamp correct cell style onto brand-new (previously empty)
# cells by copying format from a stable same-style neighbor, BEFORE
# writing their final text (Copy carries value+format; the value is
# overwritten in step 3 below).
#   style index 2 source -> B5 (plain data style, untouched by this edit)
#   style index 3 source -> E4 (highlighted style, untouched by this edit)
# ===================================================================
$ws.Range("B5").Copy($ws.Range("C5"))
$ws.Range("B5").Copy($ws.Range("A6"))
$ws.Range("B5").Copy($ws.Range("A7"))
$ws.Range("B5").Copy($ws.Range("B8"))
$ws.Range("B5").Copy($ws.Range("B9"))
$ws.Range("B5").Copy($ws.Range("C9"))
$ws.Range("B5").Copy($ws.Range("D9"))
$ws.Range("E4").Copy($ws.Range("E9"))
$ws.Range("B5").Copy($ws.Range("A10"))
$ws.Range("B5").Copy($ws.Range("B10"))
$ws.Range("B5").Copy($ws.Range("C10"))
$ws.Range("B5").Copy($ws.Range("D10"))
$ws.Range("B5").Copy($ws.Range("E10"))
$ws.Range("B5").Copy($ws.Range("A11"))
$ws.Range("B5").Copy($ws.Range("B11"))
$ws.Range("B5").Copy($ws.Range("C11"))
$ws.Range("B5").Copy($ws.Range("E11"))
$ws.Range("B5").Copy($ws.Range("B12"))
$ws.Range("B5").Copy($ws.Range("C12"))
$ws.Range("B5").Copy($ws.Range("D12"))
$ws.Range("B5").Copy($ws.Range("E12"))
$ws.Range("B5").Copy($ws.Range("A13"))
$ws.Range("B5").Copy($ws.Range("B13"))
$ws.Range("B5").Copy($ws.Range("C13"))
$ws.Range("B5").Copy($ws.Range("D13"))
$ws.Range("B5").Copy($ws.Range("E13"))

# ===================================================================
# Step 2: re-style existing cell E8 (goes from the "highlighted" style
# to the plain style)
# ===================================================================
$ws.Range("B5").Copy($ws.Range("E8"))

# ===================================================================
# Step 3: write final text - only for cells whose text actually differs
# from what is already there (re-writing unchanged multi-line cells
# would spuriously trigger Excel auto row-height).
# ===================================================================
$ws.Range("A2").Value = "text input field`nbutton (login)`nicon`ntext element"
$ws.Range("E3").Value = "login page(google)/sign up page"
$ws.Range("A4").Value = "image/profile avatar`ncashkick buttons and learn more"
$ws.Range("C4").Value = "header , main , content area , sidebar`nside navigation bar , alerts(cashkick)"
$ws.Range("C5").Value = "acceleration table"
$ws.Range("A6").Value = "buttons for reset , back navigation `nand credit review"
$ws.Range("A7").Value = "button (cancel)"
$ws.Range("B8").Value = "logout , settings ,payments"
$ws.Range("C8").Value = "popup(avatar)"
$ws.Range("D8").Value = "dashboard layout"
$ws.Range("E8").Value = "dashboard page"
$ws.Range("B9").Value = "review element"
$ws.Range("C9").Value = "success dialog"
$ws.Range("D9").Value = "success dialog layout"
$ws.Range("E9").Value = "success dialog page"
$ws.Range("A10").Value = "reset button"
$ws.Range("B10").Value = "form for forgot password"
$ws.Range("C10").Value = "reset password"
$ws.Range("D10").Value = "forgot password layout"
$ws.Range("E10").Value = "login page(password forgot)"
$ws.Range("A11").Value = "continue button"
$ws.Range("B11").Value = "reset email item"
$ws.Range("C11").Value = "continue"
$ws.Range("E11").Value = "login page(password forgot)"
$ws.Range("B12").Value = "change password form"
$ws.Range("C12").Value = "change password"
$ws.Range("D12").Value = "change password layout"
$ws.Range("E12").Value = "login page(password forgot)"
$ws.Range("A13").Value = "login button"
$ws.Range("B13").Value = "password reset item"
$ws.Range("C13").Value = "login "
$ws.Range("D13").Value = "forgot password layout"
$ws.Range("E13").Value = "login page(password forgot)"

# ===================================================================
# Step 4: rows whose multi-line cells were (re)written above picked up
# an incidental explicit row height from the write; AutoFit restores
# them to the implicit/default height used in the target layout.
# (Row 4 is intentionally NOT auto-fit: it already carries an explicit
# custom height of 51.75 in both the before and after state.)
# ===================================================================
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(6).AutoFit()

# ===================================================================
# Step 5: column width adjustments
#   col A: 20.63 -> 25.5   (exactly representable on the COM 1/6-char grid)
#   col E: 22.88 -> 25.38  (not representable exactly through ColumnWidth,
#          which this engine quantizes to 1/6-character steps; 25.333333
#          is the nearest achievable value)
# ===================================================================
$ws.Columns.Item(1).ColumnWidth = 24.65
$ws.Columns.Item(5).ColumnWidth = 24.5
